# Reshape the wide year-by-institution table (A1:O5) into a long
# "Instituição / Data / Valor" table (A1:C57), matching the upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old per-year columns (B:O held one column per "trim" period);
# columns D:O are no longer needed once the data is unpivoted into B/C.
$ws.Columns("D:O").Delete()

# New header row
$ws.Range("A1").Value = "Instituição"
$ws.Range("B1").Value = "Data"
$ws.Range("C1").Value = "Valor"

# Build the unpivoted rows: one row per (institution, period) pair,
# grouped by period (matches the order found in the workbook).
$rowCount = 56
$data = New-Object "object[,]" $rowCount,3
$data[0,0] = "Estado de Sergipe e municípios"
$data[0,1] = "2010 4° trim"
$data[0,2] = 3046.74025029022
$data[1,0] = "Governo do Estado de Sergipe"
$data[1,1] = "2010 4° trim"
$data[1,2] = 3216.0580107575706
$data[2,0] = "Município de Aracaju"
$data[2,1] = "2010 4° trim"
$data[2,2] = 22.860502551714887
$data[3,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[3,1] = "2010 4° trim"
$data[3,2] = -192.19910212804595
$data[4,0] = "Estado de Sergipe e municípios"
$data[4,1] = "2011 4° trim"
$data[4,2] = 3758.8068820088574
$data[5,0] = "Governo do Estado de Sergipe"
$data[5,1] = "2011 4° trim"
$data[5,2] = 3845.741385773845
$data[6,0] = "Município de Aracaju"
$data[6,1] = "2011 4° trim"
$data[6,2] = 12.522638399638053
$data[7,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[7,1] = "2011 4° trim"
$data[7,2] = -99.45714216462534
$data[8,0] = "Estado de Sergipe e municípios"
$data[8,1] = "2012 4° trim"
$data[8,2] = 4770.444357466845
$data[9,0] = "Governo do Estado de Sergipe"
$data[9,1] = "2012 4° trim"
$data[9,2] = 4694.092105117062
$data[10,0] = "Município de Aracaju"
$data[10,1] = "2012 4° trim"
$data[10,2] = 73.50521921131676
$data[11,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[11,1] = "2012 4° trim"
$data[11,2] = 2.847033138466494
$data[12,0] = "Estado de Sergipe e municípios"
$data[12,1] = "2013 4° trim"
$data[12,2] = 4520.2735552590975
$data[13,0] = "Governo do Estado de Sergipe"
$data[13,1] = "2013 4° trim"
$data[13,2] = 4521.617627948914
$data[14,0] = "Município de Aracaju"
$data[14,1] = "2013 4° trim"
$data[14,2] = 4.032218069450307
$data[15,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[15,1] = "2013 4° trim"
$data[15,2] = -5.376290759267076
$data[16,0] = "Estado de Sergipe e municípios"
$data[16,1] = "2014 4° trim"
$data[16,2] = 4794.213858113338
$data[17,0] = "Governo do Estado de Sergipe"
$data[17,1] = "2014 4° trim"
$data[17,2] = 4593.0140956831065
$data[18,0] = "Município de Aracaju"
$data[18,1] = "2014 4° trim"
$data[18,2] = 170.06424273250818
$data[19,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[19,1] = "2014 4° trim"
$data[19,2] = 31.135519697723595
$data[20,0] = "Estado de Sergipe e municípios"
$data[20,1] = "2015 4° trim"
$data[20,2] = 5949.231414791788
$data[21,0] = "Governo do Estado de Sergipe"
$data[21,1] = "2015 4° trim"
$data[21,2] = 5724.138661346003
$data[22,0] = "Município de Aracaju"
$data[22,1] = "2015 4° trim"
$data[22,2] = 225.4188459372783
$data[23,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[23,1] = "2015 4° trim"
$data[23,2] = -0.3260924914926435
$data[24,0] = "Estado de Sergipe e municípios"
$data[24,1] = "2016 4° trim"
$data[24,2] = 7319.488235860717
$data[25,0] = "Governo do Estado de Sergipe"
$data[25,1] = "2016 4° trim"
$data[25,2] = 7218.4810162908125
$data[26,0] = "Município de Aracaju"
$data[26,1] = "2016 4° trim"
$data[26,2] = 153.27493445986988
$data[27,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[27,1] = "2016 4° trim"
$data[27,2] = -52.28166038486509
$data[28,0] = "Estado de Sergipe e municípios"
$data[28,1] = "2017 4° trim"
$data[28,2] = 6491.191061454788
$data[29,0] = "Governo do Estado de Sergipe"
$data[29,1] = "2017 4° trim"
$data[29,2] = 6167.4090620080315
$data[30,0] = "Município de Aracaju"
$data[30,1] = "2017 4° trim"
$data[30,2] = 323.34852007338617
$data[31,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[31,1] = "2017 4° trim"
$data[31,2] = 0.41993314295244966
$data[32,0] = "Estado de Sergipe e municípios"
$data[32,1] = "2018 4° trim"
$data[32,2] = 6910.490858054237
$data[33,0] = "Governo do Estado de Sergipe"
$data[33,1] = "2018 4° trim"
$data[33,2] = 6754.901694797293
$data[34,0] = "Município de Aracaju"
$data[34,1] = "2018 4° trim"
$data[34,2] = 232.7570010253677
$data[35,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[35,1] = "2018 4° trim"
$data[35,2] = -77.16783776842384
$data[36,0] = "Estado de Sergipe e municípios"
$data[36,1] = "2019 4° trim"
$data[36,2] = 4358.394458249144
$data[37,0] = "Governo do Estado de Sergipe"
$data[37,1] = "2019 4° trim"
$data[37,2] = 4488.945876979467
$data[38,0] = "Município de Aracaju"
$data[38,1] = "2019 4° trim"
$data[38,2] = -42.36130031483485
$data[39,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[39,1] = "2019 4° trim"
$data[39,2] = -88.17760030073778
$data[40,0] = "Estado de Sergipe e municípios"
$data[40,1] = "2020 4° trim"
$data[40,2] = 6189.060050911863
$data[41,0] = "Governo do Estado de Sergipe"
$data[41,1] = "2020 4° trim"
$data[41,2] = 6229.530524458735
$data[42,0] = "Município de Aracaju"
$data[42,1] = "2020 4° trim"
$data[42,2] = 125.37582630620132
$data[43,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[43,1] = "2020 4° trim"
$data[43,2] = -165.84629985307313
$data[44,0] = "Estado de Sergipe e municípios"
$data[44,1] = "2021 4° trim"
$data[44,2] = 5935.511056055188
$data[45,0] = "Governo do Estado de Sergipe"
$data[45,1] = "2021 4° trim"
$data[45,2] = 5791.866014519512
$data[46,0] = "Município de Aracaju"
$data[46,1] = "2021 4° trim"
$data[46,2] = 137.7142424722714
$data[47,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[47,1] = "2021 4° trim"
$data[47,2] = 5.9416812635211524
$data[48,0] = "Estado de Sergipe e municípios"
$data[48,1] = "2022 4° trim"
$data[48,2] = 6082.782931423567
$data[49,0] = "Governo do Estado de Sergipe"
$data[49,1] = "2022 4° trim"
$data[49,2] = 5695.638856657847
$data[50,0] = "Município de Aracaju"
$data[50,1] = "2022 4° trim"
$data[50,2] = 262.09429603233815
$data[51,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[51,1] = "2022 4° trim"
$data[51,2] = 125.03949165056403
$data[52,0] = "Estado de Sergipe e municípios"
$data[52,1] = "2023 2° trim"
$data[52,2] = 6064.42
$data[53,0] = "Governo do Estado de Sergipe"
$data[53,1] = "2023 2° trim"
$data[53,2] = 5434.68
$data[54,0] = "Município de Aracaju"
$data[54,1] = "2023 2° trim"
$data[54,2] = 338.72
$data[55,0] = "Municípios do  Estado de Sergipe, exceto a capital"
$data[55,1] = "2023 2° trim"
$data[55,2] = 291.02

$ws.Range("A2:C" + ($rowCount + 1)).Value = $data

# Numeric formatting for the value column (2 decimal places)
$ws.Range("C2:C" + ($rowCount + 1)).NumberFormat = "0.00"

# Resize columns to fit their new (narrower/wider) contents
$ws.Columns("A:A").ColumnWidth = 44.833333333333336
$ws.Columns("B:B").ColumnWidth = 10.5
$ws.Columns("C:C").ColumnWidth = 11.833333333333334

$ws.Range("A1").Select()
